# Auto-generated Excel COM-interop edit script
# Applies header renames, title-cases connector words (de/del/la/las/los/el/y)
# in state/municipality names, fixes a float rounding cell, and trims trailing
# footer/notes rows (1202-1207), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header cells (row 1) to short English field names
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# 2) Title-case the connector words ('de', 'del', 'la', 'las', 'los', 'el', 'y')
#    in the Estado/Municipio (columns A/B) text values.
$oldNames = @(
    'Amatenango de la Frontera',
    'Amatenango del Valle',
    'Comitán de Domínguez',
    'Mazapa de Madero',
    'Ocozocoautla de Espinosa',
    'Hidalgo del Parral',
    'Ciudad de México',
    'Cuajimalpa de Morelos',
    'Nombre de Dios',
    'Pánuco de Coronado',
    'Estado de México',
    'Acambay de Ruíz Castañeda',
    'Almoloya de Alquisiras',
    'Almoloya de Juárez',
    'Atizapán de Zaragoza',
    'Chapa de Mota',
    'Coacalco de Berriozábal',
    'Ecatepec de Morelos',
    'Ixtapan de la Sal',
    'Naucalpan de Juárez',
    'San Felipe del Progreso',
    'San Martín de las Pirámides',
    'San Simón de Guerrero',
    'Tenango del Valle',
    'Tlalnepantla de Baz',
    'Valle de Bravo',
    'Valle de Chalco Solidaridad',
    'Villa de Allende',
    'Villa del Carbón',
    'San Miguel de Allende',
    'Apaseo el Grande',
    'Dolores Hidalgo Cuna de la Independencia Nacional',
    'San Francisco del Rincón',
    'San Luis de la Paz',
    'Santa Cruz de Juventino Rosas',
    'Silao de la Victoria',
    'Valle de Santiago',
    'Acapulco de Juárez',
    'Ajuchitlán del Progreso',
    'Alcozauca de Guerrero',
    'Atenango del Río',
    'Atlamajalcingo del Monte',
    'Atoyac de Álvarez',
    'Ayutla de los Libres',
    'Chilapa de Álvarez',
    'Chilpancingo de los Bravo',
    'Coahuayutla de José María Izazaga',
    'Coyuca de Benítez',
    'Coyuca de Catalán',
    'Cuetzala del Progreso',
    'Cutzamala de Pinzón',
    'Huitzuco de los Figueroa',
    'Iguala de la Independencia',
    'Ixcateopan de Cuauhtémoc',
    'Zihuatanejo de Azueta',
    'La Unión de Isidoro Montes de Oca',
    'Mártir de Cuilapan',
    'Taxco de Alarcón',
    'Técpan de Galeana',
    'Tepecoacuilco de Trujano',
    'Tlalixtaquilla de Maldonado',
    'Tlapa de Comonfort',
    'Atotonilco el Grande',
    'Cuautepec de Hinojosa',
    'Huasca de Ocampo',
    'Huejutla de Reyes',
    'Jacala de Ledezma',
    'Mineral del Chico',
    'Mineral del Monte',
    'Mixquiahuala de Juárez',
    'Omitlán de Juárez',
    'Pachuca de Soto',
    'Progreso de Obregón',
    'Santiago Tulantepec de Lugo Guerrero',
    'Tenango de Doria',
    'Tepehuacán de Guerrero',
    'Tezontepec de Aldama',
    'Tula de Allende',
    'Tulancingo de Bravo',
    'Villa de Tezontepec',
    'Zacualtipán de Ángeles',
    'Encarnación de Díaz',
    'Ixtlahuacán del Río',
    'Jilotlán de los Dolores',
    'Lagos de Moreno',
    'San Cristóbal de la Barranca',
    'San Juan de los Lagos',
    'San Miguel el Alto',
    'Talpa de Allende',
    'Tepatitlán de Morelos',
    'Tizapán el Alto',
    'Tlajomulco de Zúñiga',
    'Zacoalco de Torres',
    'Zapotlán el Grande',
    'Coalcomán de Vázquez Pallares',
    'Coatlán del Río',
    'Jonacatepec de Leandro Valle',
    'Puente de Ixtla',
    'Tetela del Volcán',
    'Tlaltizapán de Zapata',
    'Zacualpan de Amilpas',
    'Amatlán de Cañas',
    'Bahía de Banderas',
    'Acatlán de Pérez Figueroa',
    'Ayoquezco de Aldama',
    'Capulálpam de Méndez',
    'Chalcatongo de Hidalgo',
    'Coicoyán de las Flores',
    'Constancia del Rosario',
    'Cuilápam de Guerrero',
    'Fresnillo de Trujano',
    'Guevea de Humboldt',
    'Heroica Ciudad de Ejutla de Crespo',
    'Heroica Ciudad de Huajuapan de León',
    'Heroica Ciudad de Tlaxiaco',
    'Ixtlán de Juárez',
    'Heroica Ciudad de Juchitán de Zaragoza',
    'Mariscala de Juárez',
    'Mártires de Tacubaya',
    'Miahuatlán de Porfirio Díaz',
    'Mixistlán de la Reforma',
    'Nejapa de Madero',
    'Oaxaca de Juárez',
    'Ocotlán de Morelos',
    'Putla Villa de Guerrero',
    'Reforma de Pineda',
    'San Felipe Jalapa de Díaz',
    'San José del Peñasco',
    'San José del Progreso',
    'San Juan Bautista Lo de Soto',
    'San Juan del Estado',
    'San Miguel del Puerto',
    'San Miguel el Grande',
    'San Pablo Villa de Mitla',
    'Santa Cruz Tacache de Mina',
    'Santo Domingo de Morelos',
    'Tamazulápam del Espíritu Santo',
    'Tataltepec de Valdés',
    'Teococuilco de Marcos Pérez',
    'Tepelmeme Villa de Morelos',
    'Heroica Villa Tezoatlán de Segura y Luna, Cuna de la Independencia de Oaxaca',
    'Tlacolula de Matamoros',
    'Totontepec Villa de Morelos',
    'Villa de Etla',
    'Villa de Tututepec',
    'Villa de Zaachila',
    'Villa Sola de Vega',
    'Villa Talea de Castro',
    'Zimatlán de Álvarez',
    'Chalchicomula de Sesma',
    'Chila de la Sal',
    'Cuayuca de Andrade',
    'Cuetzalan del Progreso',
    'Huehuetlán el Chico',
    'Huehuetlán el Grande',
    'Huitzilan de Serdán',
    'Ixcamilpa de Guerrero',
    'Izúcar de Matamoros',
    'Los Reyes de Juárez',
    'Palmar de Bravo',
    'San Diego la Mesa Tochimiltzingo',
    'San Nicolás de los Ranchos',
    'San Salvador el Seco',
    'San Salvador el Verde',
    'Tecali de Herrera',
    'Tepanco de López',
    'Tepango de Rodríguez',
    'Tepatlaxco de Hidalgo',
    'Tepexi de Rodríguez',
    'Tepeyahualco de Cuauhtémoc',
    'Tetela de Ocampo',
    'Tlacotepec de Benito Juárez',
    'Xayacatlán de Bravo',
    'Xochitlán de Vicente Suárez',
    'Cadereyta de Montes',
    'Jalpan de Serra',
    'Landa de Matamoros',
    'Pinal de Amoles',
    'San Juan del Río',
    'Villa de Ramos',
    'Nacozari de García',
    'Contla de Juan Cuamatzi',
    'Ixtacuixtla de Mariano Matamoros',
    'Mazatecochco de José María Morelos',
    'Nanacamilpa de Mariano Arista',
    'Papalotla de Xicohténcatl',
    'San Pablo del Monte',
    'Tepetitla de Lardizábal',
    'Tetla de la Solidaridad',
    'Alto Lucero de Gutiérrez Barrios',
    'Amatlán de los Reyes',
    'Boca del Río',
    'Castillo de Teayo',
    'Cazones de Herrera',
    'Cosamaloapan de Carpio',
    'Cosautlán de Carvajal',
    'Hueyapan de Ocampo',
    'Ignacio de la Llave',
    'Ixhuatlán de Madero',
    'Ixhuatlán del Café',
    'Ixhuatlán del Sureste',
    'Juchique de Ferrer',
    'Las Vigas de Ramírez',
    'Lerdo de Tejada',
    'Martínez de la Torre',
    'Paso de Ovejas',
    'Paso del Macho',
    'Poza Rica de Hidalgo',
    'Sayula de Alemán',
    'Vega de Alatorre',
    'Zontecomatlán de López y Fuentes',
    'Zozocolco de Hidalgo'
)
$newNames = @(
    'Amatenango De La Frontera',
    'Amatenango Del Valle',
    'Comitán De Domínguez',
    'Mazapa De Madero',
    'Ocozocoautla De Espinosa',
    'Hidalgo Del Parral',
    'Ciudad De México',
    'Cuajimalpa De Morelos',
    'Nombre De Dios',
    'Pánuco De Coronado',
    'Estado De México',
    'Acambay De Ruíz Castañeda',
    'Almoloya De Alquisiras',
    'Almoloya De Juárez',
    'Atizapán De Zaragoza',
    'Chapa De Mota',
    'Coacalco De Berriozábal',
    'Ecatepec De Morelos',
    'Ixtapan De La Sal',
    'Naucalpan De Juárez',
    'San Felipe Del Progreso',
    'San Martín De Las Pirámides',
    'San Simón De Guerrero',
    'Tenango Del Valle',
    'Tlalnepantla De Baz',
    'Valle De Bravo',
    'Valle De Chalco Solidaridad',
    'Villa De Allende',
    'Villa Del Carbón',
    'San Miguel De Allende',
    'Apaseo El Grande',
    'Dolores Hidalgo Cuna De La Independencia Nacional',
    'San Francisco Del Rincón',
    'San Luis De La Paz',
    'Santa Cruz De Juventino Rosas',
    'Silao De La Victoria',
    'Valle De Santiago',
    'Acapulco De Juárez',
    'Ajuchitlán Del Progreso',
    'Alcozauca De Guerrero',
    'Atenango Del Río',
    'Atlamajalcingo Del Monte',
    'Atoyac De Álvarez',
    'Ayutla De Los Libres',
    'Chilapa De Álvarez',
    'Chilpancingo De Los Bravo',
    'Coahuayutla De José María Izazaga',
    'Coyuca De Benítez',
    'Coyuca De Catalán',
    'Cuetzala Del Progreso',
    'Cutzamala De Pinzón',
    'Huitzuco De Los Figueroa',
    'Iguala De La Independencia',
    'Ixcateopan De Cuauhtémoc',
    'Zihuatanejo De Azueta',
    'La Unión De Isidoro Montes De Oca',
    'Mártir De Cuilapan',
    'Taxco De Alarcón',
    'Técpan De Galeana',
    'Tepecoacuilco De Trujano',
    'Tlalixtaquilla De Maldonado',
    'Tlapa De Comonfort',
    'Atotonilco El Grande',
    'Cuautepec De Hinojosa',
    'Huasca De Ocampo',
    'Huejutla De Reyes',
    'Jacala De Ledezma',
    'Mineral Del Chico',
    'Mineral Del Monte',
    'Mixquiahuala De Juárez',
    'Omitlán De Juárez',
    'Pachuca De Soto',
    'Progreso De Obregón',
    'Santiago Tulantepec De Lugo Guerrero',
    'Tenango De Doria',
    'Tepehuacán De Guerrero',
    'Tezontepec De Aldama',
    'Tula De Allende',
    'Tulancingo De Bravo',
    'Villa De Tezontepec',
    'Zacualtipán De Ángeles',
    'Encarnación De Díaz',
    'Ixtlahuacán Del Río',
    'Jilotlán De Los Dolores',
    'Lagos De Moreno',
    'San Cristóbal De La Barranca',
    'San Juan De Los Lagos',
    'San Miguel El Alto',
    'Talpa De Allende',
    'Tepatitlán De Morelos',
    'Tizapán El Alto',
    'Tlajomulco De Zúñiga',
    'Zacoalco De Torres',
    'Zapotlán El Grande',
    'Coalcomán De Vázquez Pallares',
    'Coatlán Del Río',
    'Jonacatepec De Leandro Valle',
    'Puente De Ixtla',
    'Tetela Del Volcán',
    'Tlaltizapán De Zapata',
    'Zacualpan De Amilpas',
    'Amatlán De Cañas',
    'Bahía De Banderas',
    'Acatlán De Pérez Figueroa',
    'Ayoquezco De Aldama',
    'Capulálpam De Méndez',
    'Chalcatongo De Hidalgo',
    'Coicoyán De Las Flores',
    'Constancia Del Rosario',
    'Cuilápam De Guerrero',
    'Fresnillo De Trujano',
    'Guevea De Humboldt',
    'Heroica Ciudad De Ejutla De Crespo',
    'Heroica Ciudad De Huajuapan De León',
    'Heroica Ciudad De Tlaxiaco',
    'Ixtlán De Juárez',
    'Heroica Ciudad De Juchitán De Zaragoza',
    'Mariscala De Juárez',
    'Mártires De Tacubaya',
    'Miahuatlán De Porfirio Díaz',
    'Mixistlán De La Reforma',
    'Nejapa De Madero',
    'Oaxaca De Juárez',
    'Ocotlán De Morelos',
    'Putla Villa De Guerrero',
    'Reforma De Pineda',
    'San Felipe Jalapa De Díaz',
    'San José Del Peñasco',
    'San José Del Progreso',
    'San Juan Bautista Lo De Soto',
    'San Juan Del Estado',
    'San Miguel Del Puerto',
    'San Miguel El Grande',
    'San Pablo Villa De Mitla',
    'Santa Cruz Tacache De Mina',
    'Santo Domingo De Morelos',
    'Tamazulápam Del Espíritu Santo',
    'Tataltepec De Valdés',
    'Teococuilco De Marcos Pérez',
    'Tepelmeme Villa De Morelos',
    'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca',
    'Tlacolula De Matamoros',
    'Totontepec Villa De Morelos',
    'Villa De Etla',
    'Villa De Tututepec',
    'Villa De Zaachila',
    'Villa Sola De Vega',
    'Villa Talea De Castro',
    'Zimatlán De Álvarez',
    'Chalchicomula De Sesma',
    'Chila De La Sal',
    'Cuayuca De Andrade',
    'Cuetzalan Del Progreso',
    'Huehuetlán El Chico',
    'Huehuetlán El Grande',
    'Huitzilan De Serdán',
    'Ixcamilpa De Guerrero',
    'Izúcar De Matamoros',
    'Los Reyes De Juárez',
    'Palmar De Bravo',
    'San Diego La Mesa Tochimiltzingo',
    'San Nicolás De Los Ranchos',
    'San Salvador El Seco',
    'San Salvador El Verde',
    'Tecali De Herrera',
    'Tepanco De López',
    'Tepango De Rodríguez',
    'Tepatlaxco De Hidalgo',
    'Tepexi De Rodríguez',
    'Tepeyahualco De Cuauhtémoc',
    'Tetela De Ocampo',
    'Tlacotepec De Benito Juárez',
    'Xayacatlán De Bravo',
    'Xochitlán De Vicente Suárez',
    'Cadereyta De Montes',
    'Jalpan De Serra',
    'Landa De Matamoros',
    'Pinal De Amoles',
    'San Juan Del Río',
    'Villa De Ramos',
    'Nacozari De García',
    'Contla De Juan Cuamatzi',
    'Ixtacuixtla De Mariano Matamoros',
    'Mazatecochco De José María Morelos',
    'Nanacamilpa De Mariano Arista',
    'Papalotla De Xicohténcatl',
    'San Pablo Del Monte',
    'Tepetitla De Lardizábal',
    'Tetla De La Solidaridad',
    'Alto Lucero De Gutiérrez Barrios',
    'Amatlán De Los Reyes',
    'Boca Del Río',
    'Castillo De Teayo',
    'Cazones De Herrera',
    'Cosamaloapan De Carpio',
    'Cosautlán De Carvajal',
    'Hueyapan De Ocampo',
    'Ignacio De La Llave',
    'Ixhuatlán De Madero',
    'Ixhuatlán Del Café',
    'Ixhuatlán Del Sureste',
    'Juchique De Ferrer',
    'Las Vigas De Ramírez',
    'Lerdo De Tejada',
    'Martínez De La Torre',
    'Paso De Ovejas',
    'Paso Del Macho',
    'Poza Rica De Hidalgo',
    'Sayula De Alemán',
    'Vega De Alatorre',
    'Zontecomatlán De López Y Fuentes',
    'Zozocolco De Hidalgo'
)

$renameMap = @{}
for ($i = 0; $i -lt $oldNames.Length; $i++) {
    $renameMap[$oldNames[$i]] = $newNames[$i]
}

$lastDataRow = 1201
for ($r = 2; $r -le $lastDataRow; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null -and $renameMap.ContainsKey($v)) {
            $cell.Value = $renameMap[$v]
        }
    }
}

# 3) Correct the tiny floating point value in D888 (148 / 15248 recomputation)
$ws.Range("D888").Value = 0.009706190975865689

# 4) Remove the trailing footer/notes rows (1202-1207); data now ends at row 1201
$ws.Range("A1202:D1207").EntireRow.Delete()

Write-Output "edit complete"
